$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A new exposure-site row is being inserted at row 2 ("16 Bulla Rd, Essendon VIC 3040"),
# which pushes every existing data row (previously 2..159) down by one (to 3..160).
# We shift the existing rows manually (working from the bottom row upward) instead of
# using Rows.Insert(), since Insert() causes the new row to inherit the header row's
# bold/bordered style, which is not part of the target formatting.
for ($r = 159; $r -ge 2; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dest, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dest, 3).Value2 = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dest, 4).Value2 = $ws.Cells.Item($r, 4).Value2
}

# Write the new exposure-site record into the now-vacated row 2.
$ws.Cells.Item(2, 1).Value2 = "16 Bulla Rd, Essendon VIC 3040"
$ws.Cells.Item(2, 2).Value2 = -37.741667
$ws.Cells.Item(2, 3).Value2 = 144.908954
$ws.Cells.Item(2, 4).Value2 = "Moonee Valley (C)"
